$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting from H1 (bold,
# centered, bordered style) then set the text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new I/J columns, row by row (row => I,J)
$data = @{
    2  = @(1, 7)
    3  = @(1, 6)
    4  = @(1, 6)
    5  = @(1, 8)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 7)
    9  = @(1, 6)
    10 = @(1, 3)
    11 = @(1, 7)
    12 = @(1, 5)
    13 = @(1, 7)
    14 = @(1, 4)
    15 = @(1, 7)
    16 = @(1, 6)
    17 = @(1, 4)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 4)
    21 = @(1, 6)
    22 = @(7, 8)
    23 = @(1, 3)
    24 = @(3, 4)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
